# Implemented RetryAnalyzer Listener for retrying failed tests
# This script updates the Excel test-data workbook:
#  - Adds two new test rows (testOne / testTwo) on the "Tests" sheet
#  - Flips the "execute" flag for loginTestWithInValidCredentials to "no"
#  - Adjusts the "execute" flag pattern on the "DataProviderTests" sheet
#  - Updates the selected cell / dimension bookkeeping to match

$wb = $excel.ActiveWorkbook

# ---- Sheet: Tests ----
$wsTests = $wb.Worksheets.Item("Tests")

# Row 3: loginTestWithInValidCredentials execute flag changes from "yes" to "no"
$wsTests.Range("C3").Value = "no"

# New rows 4 & 5: testOne / testTwo (set column-by-column so new shared
# strings are interned in the same order Excel produced them: testOne,
# testTwo, Google search Test)
$wsTests.Range("A4").Value = "testOne"
$wsTests.Range("A5").Value = "testTwo"
$wsTests.Range("B4").Value = "Google search Test"
$wsTests.Range("B5").Value = "Google search Test"
$wsTests.Range("C4").Value = "no"
$wsTests.Range("C5").Value = "no"
$wsTests.Range("D4").Value = 0
$wsTests.Range("D5").Value = 0
$wsTests.Range("E4").Value = 1
$wsTests.Range("E5").Value = 1

$wsTests.Range("A4:E5").HorizontalAlignment = -4108
$wsTests.Range("A4:E5").VerticalAlignment = -4108

$wsTests.Range("D8").Select()

# ---- Sheet: DataProviderTests ----
$wsData = $wb.Worksheets.Item("DataProviderTests")

$wsData.Range("B2").Value = "no"
$wsData.Range("B3").Value = "yes"
$wsData.Range("B4").Value = "no"
$wsData.Range("B5").Value = "no"
$wsData.Range("B7").Value = "no"

$wsData.Range("F9").Select()
